$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "add"
$ws.Range("B9").Value = "city"
$ws.Range("C9").Value = "firdy"
$ws.Range("D9").Value = "lasf"
$ws.Range("E9").Value = " "

# F9/G9 look numeric ("55", "3333") but must stay text, matching the sibling
# rows in this sheet where similar zip/phone-like values are stored as text.
# Force text storage via NumberFormat, then restore the default "Normal"
# style so no stray style index is left on the cell.
$ws.Range("F9").NumberFormat = "@"
$ws.Range("F9").Value = "55"
$ws.Range("F9").Style = "Normal"

$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "3333"
$ws.Range("G9").Style = "Normal"
